$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7822614908218384
$ws.Range("B1").Value = 1.802002429962158
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.055581092834473
$ws.Range("E1").Value = 1.465521693229675
